$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.43877666666667
$ws.Range("H2").Value = 31.31633
$ws.Range("I2").Value = 0.478485024444405
$ws.Range("J2").Value = 0.4784850244444049
$ws.Range("M2").Value = 0.1457796666666667
$ws.Range("N2").Value = 0.437339
$ws.Range("Q2").Value = 1.521761382874445
$ws.Range("R2").Value = 13.69585244587
$ws.Range("S2").Value = 0.478485024444405
$ws.Range("T2").Value = 0.4784850244444049

$ws.Range("I3").Value = 0.07482651662844755
$ws.Range("J3").Value = 0.07482651662844754
$ws.Range("M3").Value = 0.1457796666666667
$ws.Range("N3").Value = 0.437339
$ws.Range("Q3").Value = 0.2379763160872222
$ws.Range("R3").Value = 2.141786844785
$ws.Range("S3").Value = 0.07482651662844755
$ws.Range("T3").Value = 0.07482651662844754

$ws.Range("G4").Value = 0.4778236666666666
$ws.Range("H4").Value = 1.433471
$ws.Range("I4").Value = 0.02190213241702797
$ws.Range("J4").Value = 0.02190213241702797
$ws.Range("M4").Value = 0.1457796666666667
$ws.Range("N4").Value = 0.437339
$ws.Range("Q4").Value = 0.06965697485211111
$ws.Range("R4").Value = 0.6269127736689999
$ws.Range("S4").Value = 0.02190213241702797
$ws.Range("T4").Value = 0.02190213241702797

$ws.Range("G5").Value = 9.267269333333333
$ws.Range("H5").Value = 27.801808
$ws.Range("I5").Value = 0.4247863265101195
$ws.Range("J5").Value = 0.4247863265101195
$ws.Range("M5").Value = 0.1457796666666667
$ws.Range("N5").Value = 0.437339
$ws.Range("Q5").Value = 1.350979434323555
$ws.Range("R5").Value = 12.158814908912
$ws.Range("S5").Value = 0.4247863265101195
$ws.Range("T5").Value = 0.4247863265101195
